# Apply the "updating QTL2 report code" edit:
#  1. Fix a typo in the 3rd sheet's name.
#  2. Correct several marker.id values (column H) that were shifted/mismatched
#     between the "QTL List RankZ Sex - cM", "QTL List RankZ Sex - Mbp", and
#     "QTL List RankZ SexGen - cM" sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheet 3 -----------------------------------------------------
$wsGenCm = $wb.Worksheets.Item(3)
$wsGenCm.Name = "QTL List RankZ SexGen - cM"

# --- 2. "QTL List RankZ Sex - cM" (sheet 1) marker.id corrections ----------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("H7").Value = "UNCHS035922"
$ws1.Range("H11").Value = "UNCHS005556"
$ws1.Range("H13").Value = "ICR1070"
$ws1.Range("H14").Value = "JAX00096955"
$ws1.Range("H15").Value = "UNC19016786"
$ws1.Range("H22").Value = "JAX00357615"
$ws1.Range("H29").Value = "JAX00066393r"
$ws1.Range("H31").Value = "UNCHS005556"
$ws1.Range("H39").Value = "UNCHS005556"
$ws1.Range("H44").Value = "UNC23135640"
$ws1.Range("H47").Value = "UNC3398043"

# --- 3. "QTL List RankZ Sex - Mbp" (sheet 2) marker.id correction ----------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("H29").Value = "JAX00066393"

# --- 4. "QTL List RankZ SexGen - cM" (sheet 3, just renamed) correction ----
$wsGenCm.Range("H10").Value = "UNCHS009832"
